$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.112.82'
$ws.Range("E2").Value = '  +3.03%  '
$ws.Range("D3").Value = '2.361.06'
$ws.Range("E3").Value = '  +1.12%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''311.15'
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").Value = '''108.15'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '''0.630'
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = '''0.607'
$ws.Range("E9").Value = '  -1.96%  '
$ws.Range("D10").Value = '''40.67'
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("D11").Value = '''0.0913'
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").Value = '''8.42'
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("D13").Value = '''0.109'
$ws.Range("E13").Value = '  +1.19%  '
$ws.Range("D14").Value = '''0.971'
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("D15").Value = '2.721.12'
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("D16").Value = '''15.18'
$ws.Range("E16").Value = '  -2.18%  '
$ws.Range("D17").Value = '2.371.82'
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").Value = '45.111.88'
$ws.Range("E18").Value = '  +3.23%  '
$ws.Range("D19").Value = '''14.26'
$ws.Range("E19").Value = '  +9.46%  '
$ws.Range("D20").Value = '''0.0000106'
$ws.Range("E20").Value = '  -1.16%  '
$ws.Range("D21").Value = '''7.16'
$ws.Range("E21").Value = '  -5.25%  '
$ws.Range("D22").Value = '''72.95'
$ws.Range("E22").Value = '  -1.89%  '
$ws.Range("D23").Value = '''3.50'
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("D24").Value = '''257.64'
$ws.Range("E24").Value = '  -3.89%  '
$ws.Range("D25").Value = '''2.30'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").Value = '''11.04'
$ws.Range("E27").Value = '  -1.05%  '
$ws.Range("D28").Value = '''7.19'
$ws.Range("E28").Value = '  -6.38%  '
$ws.Range("D29").Value = '''2.33'
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").Value = '''0.0972'
$ws.Range("E30").Value = '  +8.28%  '
$ws.Range("D31").Value = '''22.27'
$ws.Range("E31").Value = '  -1.29%  '
$ws.Range("D32").Value = '''36.99'
$ws.Range("E32").Value = '  -7.06%  '
$ws.Range("D33").Value = '''167.85'
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '''3.03'
$ws.Range("E34").Value = '  +6.05%  '
$ws.Range("D35").Value = '''0.130'
$ws.Range("E35").Value = '  -1.57%  '
$ws.Range("D36").Value = '''0.116'
$ws.Range("E36").Value = '  +0.42%  '
$ws.Range("D37").Value = '''4.66'
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("D38").Value = '''3.96'
$ws.Range("E38").Value = '  +4.15%  '
$ws.Range("D39").Value = '''0.0352'
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("D40").Value = '''2.88'
$ws.Range("E40").Value = '  -1.31%  '
$ws.Range("D41").Value = '''1.77'
$ws.Range("E41").Value = '  +2.90%  '
$ws.Range("D42").Value = '''98.82'
$ws.Range("E42").Value = '  -5.38%  '
$ws.Range("D43").Value = '1.893.39'
$ws.Range("E43").Value = '  +14.09%  '
$ws.Range("D44").Value = '''69.21'
$ws.Range("E44").Value = '  -3.54%  '
$ws.Range("D45").Value = '''0.228'
$ws.Range("E45").Value = '  -5.20%  '
$ws.Range("D46").Value = '''12.83'
$ws.Range("E46").Value = '  -6.63%  '
$ws.Range("D47").Value = '''1.00'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").Value = '''83.45'
$ws.Range("E48").Value = '  +9.51%  '
$ws.Range("D49").Value = '''5.62'
$ws.Range("E49").Value = '  +6.55%  '
$ws.Range("D50").Value = '''9.18'
$ws.Range("E50").Value = '  +2.17%  '
$ws.Range("D51").Value = '''109.64'
$ws.Range("E51").Value = '  -4.53%  '
